$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (old values 1002, 1003), shifting cells up
$ws.Range("A3:A4").Delete(-4162)

# Update A2's value to 91
$ws.Range("A2").Value = 91

# Move the selection to A2
$ws.Range("A2").Select()
